# "reform orders file + add opportunity to expand orders to see their items"
#
# Content changes on Лист1 (the active/only sheet):
#   - E1: 0 -> 8
#   - E2: 0 -> 6
#   - Active selection moves from D6 to G6 (leaving room to the right of the
#     existing table for the new "expand order items" columns)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 8
$ws.Range("E2").Value = 6

$ws.Range("G6").Select()
